$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data previously in row 39 now belongs in row 37, and the data
# previously in row 37 now belongs in row 39 (rows 37 and 39 swap content).

# --- Row 37: becomes the former row 39 values ---
$ws.Range("A37").Value = 131106649
$ws.Range("B37").Value = 91819
$ws.Range("D37").Value = "LC"
$ws.Range("E37").Value = 1205
$ws.Range("F37").Value = "Stor aspticka"
$ws.Range("G37").Value = "Phellinus populicola"
$ws.Range("H37").Value = "Niemelä"
$ws.Range("I37").NumberFormat = "@"
$ws.Range("I37").Value = "1"
$ws.Range("J37").Value = "mycel"
$ws.Range("Q37").Value = 601221
$ws.Range("R37").Value = 6959782
$ws.Range("X37").Value = "2025_0527"
$ws.Range("Z37").Value = "09:07"
$ws.Range("AB37").Value = "09:07"
$ws.Range("AC37").Value = "aspstubbe"

# --- Row 39: becomes the former row 37 values ---
$ws.Range("A39").Value = 131106646
$ws.Range("B39").Value = 79243
$ws.Range("D39").Value = "NT"
$ws.Range("E39").Value = 6425
$ws.Range("F39").Value = "Garnlav"
$ws.Range("G39").Value = "Alectoria sarmentosa"
$ws.Range("H39").Value = "(Ach.) Ach."
$ws.Range("I39").Value = ""
$ws.Range("J39").Value = ""
$ws.Range("Q39").Value = 601173
$ws.Range("R39").Value = 6959739
$ws.Range("X39").Value = "2025_0530"
$ws.Range("Z39").Value = "09:36"
$ws.Range("AB39").Value = "09:36"
$ws.Range("AC39").Value = "tall"
